# 15.1.1.1 Forest area workbook update
# - adds a 2023 column (O) mirroring the 2022 column (N)
# - turns 2022 (N4) into a footnoted header ("2022" + superscript "1"), stored as text
# - tweaks the "(in percent)" EN caption text
# - removes the stray trailing newline on the Kyrgyz title (A1)
# - adds footnote text (KY/RU/EN) in row 14, with a superscript "1" marker on RU/EN
# - re-centers (vertically) a number of ranges and resizes columns/rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column widths / row heights
# ---------------------------------------------------------------------------
$ws.Columns("A:C").ColumnWidth = 33.141183035714285
$ws.Columns("D:O").ColumnWidth = 8.711495535714286

$ws.Rows(1).RowHeight = 42.75
$ws.Rows(13).RowHeight = 33.75
$ws.Rows(14).RowHeight = 45

# ---------------------------------------------------------------------------
# 2. Text fixes on existing cells
# ---------------------------------------------------------------------------
# A1: drop the trailing newline in the Kyrgyz title
$ws.Range("A1").Value = "15.1.1.1 Өлкөнүн жалпы аянтына карата пайыздык катышта токойлордун аянты"

# C2: "(in percents)" -> "(in percent)"
$ws.Range("C2").Value = "(in percent)"

# ---------------------------------------------------------------------------
# 3. New column O (2023), mirroring column N (2022) values + formats
# ---------------------------------------------------------------------------
$ws.Range("D4").Copy()
$ws.Range("O4").PasteSpecial(-4122)

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("N7").Copy()
$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("N8").Copy()
$ws.Range("O8").PasteSpecial(-4122)
$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial(-4122)
$ws.Range("N10").Copy()
$ws.Range("O10").PasteSpecial(-4122)
$ws.Range("N11").Copy()
$ws.Range("O11").PasteSpecial(-4122)
$ws.Range("N12").Copy()
$ws.Range("O12").PasteSpecial(-4122)

$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 6.3
$ws.Range("O6").Value = 0.8
$ws.Range("O7").Value = 2.4
$ws.Range("O8").Value = 0.7
$ws.Range("O9").Value = 0.8
$ws.Range("O10").Value = 1
$ws.Range("O11").Value = 0.2
$ws.Range("O12").Value = 0.4

# ---------------------------------------------------------------------------
# 4. N4 becomes "2022" with a superscript footnote marker "1", stored as text
# ---------------------------------------------------------------------------
$n4 = $ws.Range("N4")
$n4.NumberFormat = "@"
$n4.Value = "20221"
$n4.HorizontalAlignment = -4152
$n4Sup = $n4.Characters(5, 1)
$n4Sup.Font.Name = "Times New Roman"
$n4Sup.Font.Size = 9
$n4Sup.Font.Bold = $true
$n4Sup.Font.Superscript = $true

# ---------------------------------------------------------------------------
# 5. Row 14 footnote text (KY plain, RU/EN with superscript "1" marker)
# ---------------------------------------------------------------------------
$a14 = $ws.Range("A14")
$a14.Value = "Кыргыз Республикасынын Өзгөчө кырдаалдар министрлигине караштуу Токой кызматынын 2022-жылдагы токой чарбасы маалыматы боюнча"

$ruText = "1По данным лесоустройства 2022 года Лесной службы при Министерстве чрезвычайных ситуаций КР"
$b14 = $ws.Range("B14")
$b14.Value = $ruText
$b14Sup = $b14.Characters(1, 1)
$b14Sup.Font.Name = "Times New Roman"
$b14Sup.Font.Size = 8
$b14Sup.Font.Superscript = $true
$b14Rest = $b14.Characters(2, $ruText.Length - 1)
$b14Rest.Font.Name = "Times New Roman"
$b14Rest.Font.Size = 8

$enText = "1According to forest management data for 2022 of the Forest Service under the Ministry of Emergency Situations of the Kyrgyz Republic"
$c14 = $ws.Range("C14")
$c14.Value = $enText
$c14Sup = $c14.Characters(1, 1)
$c14Sup.Font.Name = "Times New Roman"
$c14Sup.Font.Size = 8
$c14Sup.Font.Superscript = $true
$c14Rest = $c14.Characters(2, $enText.Length - 1)
$c14Rest.Font.Name = "Times New Roman"
$c14Rest.Font.Size = 8

# ---------------------------------------------------------------------------
# 6. Vertical centering touch-ups (matches the refreshed style set)
# ---------------------------------------------------------------------------
$ws.Range("D4:O4").VerticalAlignment = -4108
$ws.Range("D5:O12").VerticalAlignment = -4108
$ws.Range("A13:C13").VerticalAlignment = -4108
$ws.Range("A14:C14").VerticalAlignment = -4108

Write-Output "done"
